# Conduct S3 tests for Graphql with resolved fields
#
# Adds a "Srednia" (AVERAGE) column next to the raw trial results on every
# sheet, and fixes two cells on "1 rekord" that were mistakenly stored as
# text ("777MB" / "722MB") instead of plain numbers (777 / 722).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "1 rekord": 4 trial columns (G:J) -> average goes to column K
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("1 rekord")

# Two cells were recorded as the text "777MB"/"722MB"; they should just be
# the plain numbers so that the new AVERAGE() formulas can use them.
$ws1.Range("I37").Value = 777
$ws1.Range("J37").Value = 777
$ws1.Range("I38").Value = 722
$ws1.Range("J38").Value = 722

for ($row = 9; $row -le 38; $row++) {
    $cell = $ws1.Cells.Item($row, 11)  # column K
    $cell.Formula = "=AVERAGE(G" + $row + ":J" + $row + ")"
}

$ws1.Range("I39").Select()

# ---------------------------------------------------------------------
# Sheet "100 rekordów": 3 trial columns (G:I) -> average goes to column J
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("100 rekordów")

for ($row = 9; $row -le 38; $row++) {
    $cell = $ws2.Cells.Item($row, 10)  # column J
    $cell.Formula = "=AVERAGE(G" + $row + ":I" + $row + ")"
}

$ws2.Columns.Item(10).ColumnWidth = 13.3

$ws2.Range("L24").Select()

# ---------------------------------------------------------------------
# Sheet "500 rekordów": 3 trial columns (G:I) -> average goes to column J
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("500 rekordów")

for ($row = 9; $row -le 38; $row++) {
    $cell = $ws3.Cells.Item($row, 10)  # column J
    $cell.Formula = "=AVERAGE(G" + $row + ":I" + $row + ")"
}

$ws3.Activate()
$ws3.Range("M20").Select()
